$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header values
$ws.Range("A1").Value = "Nombre País"
$ws.Range("B1").Value = "Total de Ventas"
$ws.Range("C1").Value = "Monto de Facturación"

# Column widths (approx matching target)
$ws.Columns.Item(1).ColumnWidth = 18.5546875
$ws.Columns.Item(2).ColumnWidth = 21.21875
$ws.Columns.Item(3).ColumnWidth = 24.88671875

# Row height
$ws.Rows.Item(1).RowHeight = 19.2

# Formatting: bold, size 12 font, fill color (Accent6, Lighter 60%), vertical center
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12
$headerRange.Font.Name = "Calibri"
$headerRange.Interior.ThemeColor = 9
$headerRange.Interior.TintAndShade = 0.59999389629810485
$headerRange.VerticalAlignment = -4108
